# Corrected search component to fit home page
# Re-signs the "Amount" (column I) of every MoMo transaction row where the
# account owner (FRI:250783975669/MSISDN) is the sender, turning those
# amounts negative so outgoing money is distinguished from incoming money.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose "From" (column D) is the account owner -> Amount becomes negative.
$rowsToNegate = @(5, 8, 10, 12, 16, 17, 18, 19, 21, 23, 24, 26, 44, 46, 47, 52, 64, 65, 66, 67)

foreach ($r in $rowsToNegate) {
    $cell = $ws.Cells.Item($r, 9)   # Column I = Amount
    $current = $cell.Value2
    if ($current -gt 0) {
        $cell.Value = -$current
    }
}

# Resize column K ("Currency") to fit its contents, mirroring the sheet's
# other bestFit columns.
$ws.Columns("K").AutoFit()

# Update the saved selection / scroll position of the sheet view.
$ws.Activate()
$ws.Range("D19").Select()
